# Applies the "Automatic update of files." commit: the data rows (2..50) of the
# "Avverkningsanmälningar" sheet get re-ordered (the underlying source feed re-sorted
# the records), and every row's "Förändrad" timestamp (column C) is bumped from
# 2026-02-25 (46078) to 2026-02-26 (46079).
#
# Strategy: snapshot every existing data row (columns A..Z) into memory first
# (so overlapping row-to-row moves never read already-overwritten data), then
# write the snapshots back out in their new row positions, and finally stamp
# column C with the new date serial for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Destination row -> source row (where the record that ends up in the destination
# row currently lives).
$rowMap = @{
  2=2; 3=3; 4=4; 5=6; 6=5; 7=7; 8=14; 9=17; 10=10; 11=15; 12=11; 13=8; 14=9; 15=12;
  16=13; 17=16; 18=18; 19=19; 20=20; 21=48; 22=34; 23=49; 24=44; 25=33; 26=41; 27=21;
  28=31; 29=26; 30=32; 31=39; 32=46; 33=45; 34=23; 35=25; 36=35; 37=27; 38=28; 39=29;
  40=42; 41=22; 42=24; 43=30; 44=38; 45=37; 46=36; 47=40; 48=50; 49=47; 50=43
}

$firstRow = 2
$lastRow = 50

# 1) Snapshot every row's contents before making any changes.
$snapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowObj = @{}
    $rowObj.A = $ws.Range("A$r").Value()
    $rowObj.B = $ws.Range("B$r").Value2()
    $rowObj.D = $ws.Range("D$r").Value()
    $rowObj.E = $ws.Range("E$r").Value()
    $rowObj.HasF = -not [string]::IsNullOrEmpty($ws.Range("F$r").Value())
    $rowObj.F = $ws.Range("F$r").Value()
    for ($col = 7; $col -le 17; $col++) {
        $letter = [char](64 + $col)
        $rowObj["C$letter"] = $ws.Cells.Item($r, $col).Value()
    }
    $rowObj.R = $ws.Range("R$r").Value()
    $rowObj.HasLinks = -not [string]::IsNullOrEmpty($ws.Range("S$r").Formula)
    foreach ($letter in @("S","T","V","W","X","Y")) {
        $rowObj["F_$letter"] = $ws.Range("$letter$r").Formula
    }
    $rowObj.HasZ = -not [string]::IsNullOrEmpty($ws.Range("Z$r").Formula)
    $rowObj.F_Z = $ws.Range("Z$r").Formula
    $snapshots[$r] = $rowObj
}

# 2) Write every destination row from the matching source snapshot.
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $data = $snapshots[$srcRow]

    $ws.Range("A$destRow").Value = $data.A
    $ws.Range("B$destRow").Value = $data.B
    $ws.Range("D$destRow").Value = $data.D
    $ws.Range("E$destRow").Value = $data.E

    if ($data.HasF) {
        $ws.Range("F$destRow").Value = $data.F
    } else {
        $ws.Range("F$destRow").ClearContents()
    }

    for ($col = 7; $col -le 17; $col++) {
        $letter = [char](64 + $col)
        $ws.Cells.Item($destRow, $col).Value = $data["C$letter"]
    }

    if (-not [string]::IsNullOrEmpty($data.R)) {
        $ws.Range("R$destRow").Value = $data.R
    }

    if ($data.HasLinks) {
        foreach ($letter in @("S","T","V","W","X","Y")) {
            $ws.Range("$letter$destRow").Formula = $data["F_$letter"]
        }
    } else {
        foreach ($letter in @("S","T","V","W","X","Y")) {
            $ws.Range("$letter$destRow").ClearContents()
        }
    }

    if ($data.HasZ) {
        $ws.Range("Z$destRow").Formula = $data.F_Z
    } else {
        $ws.Range("Z$destRow").ClearContents()
    }
}

# 3) Every row's "Förändrad" (column C) becomes the same new date serial.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("C$r").Value = 46079
}
